$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (COMP3804): populate new "Final Grade" column E ---
$ws1.Range("E3").Value = 7.6
$ws1.Range("E4").Value = 9.545454545454547
$ws1.Range("E5").Value = 10
$ws1.Range("E6").Value = 8.6
$ws1.Range("E7").Value = 12.75
$ws1.Range("E8").Value = 8.289473684210527
$ws1.Range("E9").Value = 24

# --- Sheet2 (COMP3804 Copy): update grades received in column D, then populate E ---
$ws2.Range("D3").Value = 86.56
$ws2.Range("D4").Value = 100
$ws2.Range("D5").Value = 0.84
$ws2.Range("D7").ClearContents()
$ws2.Range("D8").NumberFormat = "0%"
$ws2.Range("D8").Value = 1

$ws2.Range("E3").Value = 8.656000000000001
$ws2.Range("E4").Value = 10
$ws2.Range("E5").Value = 8.4
$ws2.Range("E6").Value = 8.6
$ws2.Range("E7").Value = 0
$ws2.Range("E8").Value = 15
$ws2.Range("E9").Value = 24

# New row 10: final grade cell, formatted as percentage, left blank
$ws2.Range("D10").NumberFormat = "0%"

# Make sheet2 the active tab/sheet, with D10 selected
$ws2.Activate()
$ws2.Range("D10").Select()
